# Update BOM with price per assembly
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column O: "Production cost" = Qty * Unit price ---
$ws.Range("O1").Value = "Production cost"

# --- Fill in the previously-empty eBay production-cost row (row 23) ---
$ws.Range("G23").Value = "(eBay)"
$ws.Range("L23").Value = 50
$ws.Range("M23").Value = 8.9
$ws.Range("N23").Value = 445

$ws.Range("O2").Formula = "=B2*M2"
$ws.Range("O3:O40").Formula = "=B3*M3"
$ws.Range("O40").NumberFormat = $ws.Range("O39").NumberFormat

# --- Grand total row ---
$ws.Range("O41").Formula = "=SUM(O2:O40)"
$ws.Range("O41").NumberFormat = $ws.Range("O39").NumberFormat
$ws.Range("O41").Font.Bold = $true
$ws.Rows.Item(41).RowHeight = 41

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 41.1666666
$ws.Columns.Item(15).ColumnWidth = 14

# --- Selection / view state ---
$ws.Range("G25").Select()
